# Update the header row of the sheet:
#  - A1's text changes from "P35" to "Valor"
#  - A1:B1 get centered horizontal alignment
#  - A1:B1 becomes the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 keeps referencing the same shared-string slot, just with new text "Valor"
$ws.Range("A1").Value = "Valor"

# Center-align the header cells (creates/reuses a cellXfs entry with horizontal="center")
$ws.Range("A1:B1").HorizontalAlignment = -4108  # xlCenter

# Make A1:B1 the current selection, matching the saved sheetView selection
$ws.Range("A1:B1").Select()
